# Updates test-result log values (Result/Date columns) across the ABP
# PaymentsCC Katalon data sheets, reflecting the latest automated test run
# (Object repository additions for "Pay As a guest" page / ABP suite).

$wb = $excel.ActiveWorkbook

function Set-ResultRow {
    param(
        [string]$SheetName,
        [string]$Cell,
        [string]$Value
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range($Cell).Value = $Value
}

# VerifySuccessfulPaymentCCNoCF
Set-ResultRow "VerifySuccessfulPaymentCCNoCF" "B2" "Thu Sep 25 13:22:50 IST 2025"
Set-ResultRow "VerifySuccessfulPaymentCCNoCF" "A3" "Pass"
Set-ResultRow "VerifySuccessfulPaymentCCNoCF" "B3" "Thu Sep 25 13:25:33 IST 2025"
Set-ResultRow "VerifySuccessfulPaymentCCNoCF" "A4" "Pass"
Set-ResultRow "VerifySuccessfulPaymentCCNoCF" "B4" "Thu Sep 25 13:28:15 IST 2025"

# VerifyCFVerbiageOnRPCCDCF
Set-ResultRow "VerifyCFVerbiageOnRPCCDCF" "B2" "Thu Sep 25 13:14:04 IST 2025"

# VerifyStaticTextOnRecieptCC
Set-ResultRow "VerifyStaticTextOnRecieptCC" "B2" "Thu Sep 25 13:19:03 IST 2025"

# VerifySuccessfulPaymentCCSCF
Set-ResultRow "VerifySuccessfulPaymentCCSCF" "B2" "Thu Sep 25 13:30:46 IST 2025"
Set-ResultRow "VerifySuccessfulPaymentCCSCF" "B3" "Thu Sep 25 13:33:20 IST 2025"
Set-ResultRow "VerifySuccessfulPaymentCCSCF" "B4" "Thu Sep 25 13:36:16 IST 2025"

# VerifySuccessfulPaymentCCDCF
Set-ResultRow "VerifySuccessfulPaymentCCDCF" "B2" "Tue Sep 30 18:23:15 IST 2025"
Set-ResultRow "VerifySuccessfulPaymentCCDCF" "B3" "Tue Sep 30 18:23:57 IST 2025"
Set-ResultRow "VerifySuccessfulPaymentCCDCF" "A4" "Pass"
Set-ResultRow "VerifySuccessfulPaymentCCDCF" "B4" "Tue Sep 30 18:24:34 IST 2025"

# VerifySuccessfulPaymentSPMCCDCF
Set-ResultRow "VerifySuccessfulPaymentSPMCCDCF" "A2" "Fail"
Set-ResultRow "VerifySuccessfulPaymentSPMCCDCF" "B2" "Thu Sep 25 13:56:22 IST 2025"

# UiVerificationSPPaymentMethodCC
Set-ResultRow "UiVerificationSPPaymentMethodCC" "B2" "Thu Sep 25 14:06:16 IST 2025"

# VerifyStaticTextOnPPCCDCF
Set-ResultRow "VerifyStaticTextOnPPCCDCF" "B2" "Thu Sep 25 13:16:53 IST 2025"

# CreateVerifyDeleteCCPM
Set-ResultRow "CreateVerifyDeleteCCPM" "A2" "Pass"
Set-ResultRow "CreateVerifyDeleteCCPM" "B2" "Thu Sep 25 13:09:57 IST 2025"
